$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the "Date" metadata value.
$ws.Range("B8").Value2 = "2024-07-01T07:50:29+00:00"

# 2. Insert a new "Jurisdiction" row (with an empty value) right after the
#    "Contact" row (row 10), pushing every following row down by one.
#    Shift rows 21..11 down to 22..12 one at a time (bottom-up so we never
#    overwrite data we still need to read). Clear the destination row first
#    since pasting a range that contains an empty cell does not clear an
#    already-populated destination cell.
for ($r = 21; $r -ge 11; $r--) {
    $ws.Range("A$($r+1):B$($r+1)").ClearContents()
    $ws.Range("A$($r):B$($r)").Copy($ws.Range("A$($r+1)"))
}

$ws.Range("A11").Value2 = "Jurisdiction"
$ws.Range("B11").Value2 = ""
